$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellB = $ws.Cells.Item($r, 2)

    $valA = $cellA.Value2
    if ($valA -ne $null) {
        $cellA.Value2 = $valA.ToUpper()
    }

    $valB = $cellB.Value2
    if ($valB -ne $null) {
        $cellB.Value2 = $valB.ToUpper()
    }
}

# Header text changes for column C (not just case change)
$ws.Cells.Item(1, 3).Value2 = "PREMIUM TOTAL"
